$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 28 (SC 92) first, then row 26 (RM 232), deleting bottom-up
# so earlier row indices remain valid after each delete.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# After deletion, data that was in rows 27-35 shifted up to 26-33 (two rows removed).
# Update E26 (was -5 for SC5) to empty, E27 (was empty for SC101) to -10,
# and E29 (was -6.8 for SC119) to empty.
$ws.Range("E26").ClearContents()
$ws.Range("E27").Value = -10
$ws.Range("E29").ClearContents()
